$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: "Interactivity in R" -> "Building Shiny applications"; link_it flips to TRUE
$ws.Range("D20").Value = "Building Shiny applications"
$ws.Range("C20").Value = $true

# Row 21: "Shiny applications" -> "Building Shiny applications (part II)"
$ws.Range("D21").Value = "Building Shiny applications (part II)"

# Move the active selection to C21 (matches the saved view state in the diff)
$ws.Range("C21").Select() | Out-Null
